$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1): new columns I ("I0") and J ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from an existing
# header cell (H1) onto the two new header cells, reusing the same cell
# style rather than creating new style records.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I and J columns, rows 2-8 (both columns hold identical
# values per the source diff)
$values = @(10, 9, 5, 6, 9, 9, 9)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
